# Applies the cryptos.xlsx price-refresh update described by the commit diff.
# Numeric-look text cells (e.g. "574.63") are routed through a scratch cell
# that is explicitly formatted as Text, then copied in via PasteSpecial(values)
# so the target cell keeps plain "General" formatting/style (no quote-prefix,
# no number coercion) exactly like the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

# Row 2
$ws.Range("D2").Value = '66.256.67'
$ws.Range("E2").Value = '  -1.13%  '

# Row 3
$ws.Range("D3").Value = '3.078.45'
$ws.Range("E3").Value = '  -1.44%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$scratch.NumberFormat = "@"
$scratch.Value = '574.63'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E5").Value = '  -0.49%  '

# Row 6
$scratch.NumberFormat = "@"
$scratch.Value = '170.09'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E6").Value = '  -1.78%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = '3.075.53'
$ws.Range("E8").Value = '  -1.35%  '

# Row 9
$ws.Range("E9").Value = '  -2.36%  '

# Row 10
$scratch.NumberFormat = "@"
$scratch.Value = '6.32'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E10").Value = '  -1.99%  '

# Row 11
$ws.Range("E11").Value = '  -2.83%  '

# Row 12
$ws.Range("E12").Value = '  -2.62%  '

# Row 13
$ws.Range("E13").Value = '  -3.25%  '

# Row 14
$ws.Range("E14").Value = '  -3.71%  '

# Row 15
$ws.Range("E15").Value = '  -1.45%  '

# Row 16
$ws.Range("D16").Value = '3.592.76'
$ws.Range("E16").Value = '  -1.29%  '

# Row 17
$ws.Range("D17").Value = '66.221.32'
$ws.Range("E17").Value = '  -1.10%  '

# Row 18
$ws.Range("E18").Value = '  -3.11%  '

# Row 19
$scratch.NumberFormat = "@"
$scratch.Value = '16.60'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E19").Value = '  +1.84%  '

# Row 20
$ws.Range("D20").Value = '3.078.50'
$ws.Range("E20").Value = '  -1.29%  '

# Row 21
$scratch.NumberFormat = "@"
$scratch.Value = '487.16'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E21").Value = '  +2.26%  '

# Row 22
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$scratch.NumberFormat = "@"
$scratch.Value = '0.686'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E22").Value = '  -3.51%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$scratch.NumberFormat = "@"
$scratch.Value = '7.68'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E23").Value = '  -3.35%  '

# Row 24
$scratch.NumberFormat = "@"
$scratch.Value = '82.52'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E24").Value = '  -1.78%  '

# Row 26
$ws.Range("E26").Value = '  -3.75%  '

# Row 27
$ws.Range("E27").Value = '  +1.40%  '

# Row 28
$ws.Range("E28").Value = '  +0.01%  '

# Row 29
$ws.Range("E29").Value = '  -0.31%  '

# Row 30
$scratch.NumberFormat = "@"
$scratch.Value = '2.24'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E30").Value = '  -5.41%  '

# Row 31
$ws.Range("E31").Value = '  -3.07%  '

# Row 32
$ws.Range("E32").Value = '  -2.98%  '

# Row 33
$ws.Range("E33").Value = '  -3.98%  '

# Row 34
$ws.Range("D34").Value = '0.0₃0899'
$ws.Range("E34").Value = '  -5.38%  '

# Row 35
$scratch.NumberFormat = "@"
$scratch.Value = '1.00'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E35").Value = '  +0.07%  '

# Row 36
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$scratch.NumberFormat = "@"
$scratch.Value = '0.948'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E36").Value = '  -2.87%  '

# Row 37
$ws.Range("B37").Value = 'Arweave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$scratch.NumberFormat = "@"
$scratch.Value = '47.17'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E37").Value = '  +0.21%  '

# Row 38
$scratch.NumberFormat = "@"
$scratch.Value = '5.56'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E38").Value = '  -5.15%  '

# Row 39
$ws.Range("E39").Value = '  -1.10%  '

# Row 40
$ws.Range("E40").Value = '  -5.03%  '

# Row 41
$scratch.NumberFormat = "@"
$scratch.Value = '0.299'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E41").Value = '  -4.34%  '

# Row 42
$ws.Range("E42").Value = '  -4.80%  '

# Row 43
$ws.Range("D43").Value = '2.785.76'
$ws.Range("E43").Value = '  -1.08%  '

# Row 44
$ws.Range("E44").Value = '  -2.80%  '

# Row 45
$ws.Range("E45").Value = '  -1.94%  '

# Row 46
$ws.Range("E46").Value = '  -0.66%  '

# Row 47
$scratch.NumberFormat = "@"
$scratch.Value = '364.23'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E47").Value = '  -4.93%  '

# Row 48
$ws.Range("E48").Value = '  +0.01%  '

# Row 49
$scratch.NumberFormat = "@"
$scratch.Value = '24.42'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E49").Value = '  -1.86%  '

# Row 50
$ws.Range("E50").Value = '  -2.26%  '

# Row 51
$ws.Range("E51").Value = '  -2.31%  '
